# Refresh the crypto price/volume table (columns D and E) for every coin row,
# and fix the Dai/LEO and FirstDigitalUSD/Filecoin row ordering (columns B-E).
# Numeric-looking text such as '579.13' is written with a leading apostrophe so
# Excel keeps it as text instead of silently converting it to a float.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.521.06"
$ws.Range("E2").Value = "  +0.82%  "
$ws.Range("D3").Value = "3.443.70"
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'579.13"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("D6").Value = "'149.67"
$ws.Range("E6").Value = "  +9.14%  "
$ws.Range("D7").Value = "3.445.02"
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("D10").Value = "'7.82"
$ws.Range("E10").Value = "  +2.72%  "
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "4.031.57"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("D14").Value = "'27.99"
$ws.Range("E14").Value = "  +6.73%  "
$ws.Range("E15").Value = "  -0.32%  "
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("D17").Value = "3.440.48"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "61.598.84"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("E19").Value = "  +8.49%  "
$ws.Range("D20").Value = "'14.30"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'388.22"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").Value = "'0.570"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("D24").Value = "3.585.27"
$ws.Range("E24").Value = "  +1.68%  "
$ws.Range("D25").Value = "'72.90"
$ws.Range("E25").Value = "  +2.28%  "
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "'5.78"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  -0.30%  "
$ws.Range("E28").Value = "  -1.46%  "
$ws.Range("D29").Value = "'0.181"
$ws.Range("E29").Value = "  +4.36%  "
$ws.Range("D30").Value = "'7.76"
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'1.52"
$ws.Range("E32").Value = "  -13.74%  "
$ws.Range("D33").Value = "'8.25"
$ws.Range("E33").Value = "  +1.14%  "
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'24.03"
$ws.Range("E36").Value = "  +1.26%  "
$ws.Range("D37").Value = "'5.31"
$ws.Range("E37").Value = "  +1.62%  "
$ws.Range("D38").Value = "'7.06"
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +1.11%  "
$ws.Range("D40").Value = "'166.42"
$ws.Range("E40").Value = "  +1.12%  "
$ws.Range("D41").Value = "'0.0792"
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("D42").Value = "'26.54"
$ws.Range("E42").Value = "  +9.65%  "
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.50"
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").Value = "2.611.09"
$ws.Range("E48").Value = "  +5.57%  "
$ws.Range("E49").Value = "  -3.07%  "
$ws.Range("D50").Value = "'7.04"
$ws.Range("E50").Value = "  +3.51%  "
$ws.Range("D51").Value = "'23.24"
$ws.Range("E51").Value = "  +0.08%  "
